$d = $word.ActiveDocument

function Collapse-Run {
    # Simple case: the run that needs merging is fully interior to the
    # paragraph (there is a run AFTER the trailing proofErr mark), so a
    # plain Find/Replace over the exact text correctly swallows the
    # <w:proofErr/> markers sitting between the runs.
    param([string]$text)
    $d.Content.Find.Execute($text, $true, $false, $false, $false, $false,
        $true, 1, $false, $text, 2) | Out-Null
}

function Collapse-Run-TrailingMark {
    # Harder case: the <w:proofErr w:type="gramEnd"/> is the LAST child of
    # the paragraph (no run follows it), so it sits just outside the Find
    # match range and a normal Find/Replace leaves it behind. Push a
    # throw-away character past the end of the match first (so the stray
    # mark becomes an interior node), then redo the Find to include that
    # throw-away character in the search text (but not in the replacement)
    # -- this absorbs and discards the orphaned proofErr mark along with it.
    param([string]$text)
    $r = $d.Content
    $r.Find.Execute($text, $true, $false, $false, $false, $false,
        $true, 1, $false, "", 0) | Out-Null
    $r.InsertAfter("X")
    $d.Content.Find.Execute($text + "X", $true, $false, $false, $false, $false,
        $true, 1, $false, $text, 2) | Out-Null
}

# --- 1. Remove stray grammar-check run splits (collapse to single runs) ---

Collapse-Run "Spawners appear in corrupted lands, the player must find them and destroy them before they grow"

Collapse-Run-TrailingMark "More enemies and tougher enemies spawn as the player progresses"

Collapse-Run ("Squads " + [char]0x2013 + " allows soldiers, archers and priests to be commanded as groups rather than as single units. Combine squads to create armies")

Collapse-Run ("Explosive " + [char]0x2013 + " moves really slowly but deals AOE damage when in range of followers or buildings. They must be killed in time")

Collapse-Run "Check if target is static, if so there is no need to update every frame"

Collapse-Run-TrailingMark "Buildings can call nearby followers to defend it?"

Collapse-Run "Move towards target unless there are corpses in range and the raise dead ability is off cooldown (if so move towards corpses)"

Collapse-Run "Rework HUD (needs to be easier to update intuitively (finds all of the data from static classes etc.)"

# --- 2. Insert new "Save corpses" bullet right before "Store status effects" ---

$r1 = $d.Content
$r1.Find.Execute("Store status effects", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r1.InsertBefore("Save corpses`r")

# --- 3. Insert new "Day/night cycle" bullet after the last bullet in the doc ---

$r2 = $d.Content
$r2.Find.Execute("Add comment headers to separate sections in classes", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r2.InsertAfter("`rDay/night cycle")
